# "Mise à jour de l'auto-évaluation"
#
# - add reviewer comments (column C) on several checklist rows
# - change an existing reviewer comment's text
# - flip three checklist status cells from "not done" (red) to "done" (green)
# - move the active selection from C37 to A7 (and let the view scroll back to top)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New remarks added in column C (ordered by row so the shared-string table
# grows in the same order the rows appear in).
$ws.Range("C3").Value = "demander si juste partie 1 ou tout le projet"
$ws.Range("C6").Value = "vanessa"
$ws.Range("C24").Value = "Raphaël"
$ws.Range("C29").Value = "Raphaël : Contrôleur, Services (non DAO)   /   Alex :  Entités JPA, DTO, Services DAO"
$ws.Range("C30").Value = "Raphael et Alex"

# Mark B4, B5 and B8 as completed: red fill -> green fill (same green already
# used elsewhere in the sheet, e.g. B12:B25).
$ws.Range("B4").Interior.Color = RGB(0, 176, 80)
$ws.Range("B5").Interior.Color = RGB(0, 176, 80)
$ws.Range("B8").Interior.Color = RGB(0, 176, 80)

# Update the selected cell / scroll position.
$ws.Range("A7").Select()
